# NP comparisons within seasons
# Applies:
#  - New significance-letter cells (I6, K6, I15, I24, K24, I35, K35) using "*","**","***","****"
#  - Removes stray cell borders on several header/body cells (fill stays, border goes)
#  - Updates the sheet's scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the new significance-marker cells.
#    Introduce the brand-new shared strings in ascending length order ("*",
#    "**", "***", "****") first so they land in the shared-string table in
#    that order, then fill in the remaining cells that reuse them.
# ---------------------------------------------------------------------------
$ws.Range("K6").Value2  = "*"
$ws.Range("I24").Value2 = "**"
$ws.Range("I35").Value2 = "***"
$ws.Range("I6").Value2  = "****"

$ws.Range("I15").Value2 = "****"
$ws.Range("K24").Value2 = "**"
$ws.Range("K35").Value2 = "**"

# ---------------------------------------------------------------------------
# 2) Strip the (invisible, borderId=0) "apply border" formatting left on a
#    handful of cells so they share the plain fill-only style.
# ---------------------------------------------------------------------------
$noBorderCells = @("G4","B5","B6","C7","B8","G13","G22","G33","B7")
foreach ($addr in $noBorderCells) {
    $ws.Range($addr).Borders.LineStyle = -4142   # xlLineStyleNone
}

# ---------------------------------------------------------------------------
# 3) Update the view: scroll so row 3 is at the top and select G6:M35.
# ---------------------------------------------------------------------------
$ws.Range("G6:M35").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
